# Apply targeted cell updates to the FRIDAY schedule sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AV2").Value = "RLN"
$ws.Range("Y4").Value = 32
$ws.Range("P5").Value = "TT"
$ws.Range("AC5").Value = 29
$ws.Range("AF5").Value = "TAC/BJM"
$ws.Range("AJ5").Value = "TT"
$ws.Range("AN5").Value = "TT"
$ws.Range("AB6").Value = "AMS/CS/SL"
$ws.Range("AF6").Value = "AMS/CS/SL"
$ws.Range("AC8").Value = 26
$ws.Range("P11").Value = "MP"
$ws.Range("AC11").Value = 41
$ws.Range("AJ11").Value = "MP"
$ws.Range("AN11").Value = "MP"
$ws.Range("Y12").Value = 40
$ws.Range("AC12").Value = 27
$ws.Range("AC14").Value = 68
$ws.Range("Y15").Value = 29
$ws.Range("AC15").Value = 63
$ws.Range("AC16").Value = 57
$ws.Range("AC17").Value = 2
$ws.Range("AC19").Value = 54
$ws.Range("AC22").Value = 1
$ws.Range("AB23").Value = "BS/LP/AS"
$ws.Range("AC23").Value = 69
$ws.Range("AC24").Value = 62
$ws.Range("AB25").Value = "FGN/CK"
$ws.Range("AF25").Value = "FGN/CK"
